$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CARTE-DES-VINS")

$ws.Columns.Item(10).Cut() | Out-Null
$ws.Columns.Item(2).Insert() | Out-Null
